$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.083.37'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.831.52'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.74'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.46%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6334'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.94%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07544'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.60%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2948'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.20%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.15'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07706'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.84%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.822.47'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.006'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6705'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.86%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '83.32'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009642'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.076'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.62%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.094.38'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.80%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.62'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '226.55'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.79%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.163'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.002'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '160.91'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.26%  '
$ws.Range('E25').Value = '  +4.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.534'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.96'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.85%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.506'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.78%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.152'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.068'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.05483'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.58%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.205'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.21%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.863'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7457'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.26%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.142'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.92%  '
$ws.Range('E36').Value = '  +1.86%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.245.96'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.10%  '
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01786'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.641'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9031'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.60%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.001'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.49'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.976.86'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.17%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '65.17'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.23%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000122'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.68%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5106'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4068'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.28%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.004'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.661'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.07%  '
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.788'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.54%  '
